$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$companies = @{
    '2' = 'Chizhou Lihua Environment Technology Co., Ltd.'
    '3' = 'Shenzhen Kingyard Int''l Trading Co., Ltd.'
    '4' = 'Shanghai Chuangsi Youxiang Trading Co., Ltd.'
    '5' = 'Henrich (shandong) Health Technology Co., Ltd.'
    '6' = 'Hebei Dirui Textile Co., Ltd.'
    '7' = 'Changshu Changfu Textile Co., Ltd.'
    '8' = 'Liaoning Tanghe Daily Necessities Co., Ltd.'
    '9' = 'Huizhou MeiJi Non Woven Fabric Products Co.,Ltd.'
    '10' = 'Shanghai Chuangsi Youxiang Trading Co., Ltd.'
    '11' = 'Hangzhou Biogen Hygiene Co., Ltd.'
    '12' = 'Shanghai Chuangsi Youxiang Trading Co., Ltd.'
    '13' = 'Chizhou Lihua Environment Technology Co., Ltd.'
    '14' = 'Ningbo Riway Nonwovens Tech Co., Ltd.'
    '15' = 'Hangzhou Biogen Hygiene Co., Ltd.'
    '16' = 'Henan Enda New Material Co., Ltd.'
    '17' = 'Zhejiang Wipex New Material Technology Co., Ltd.'
    '18' = 'Blueseth Technology (Dalian) Co., Ltd.'
    '19' = 'Shaoxing Elite Bio-Tech Co., Ltd.'
    '20' = 'Hangzhou Clear Cleaning Products Co., Ltd.'
    '21' = 'Shanghai Chuangsi Youxiang Trading Co., Ltd.'
    '22' = 'Changshu Changfu Textile Co., Ltd.'
    '23' = 'Jinhua Zhili Daily Necessities Co., Ltd.'
    '24' = 'Shijiazhuang Handu Garment Co., Ltd.'
    '25' = 'Zhejiang E-Sun Enviromental Technology Co., Ltd.'
    '26' = 'Hangzhou Wipex Nonwovens Co., Ltd.'
    '27' = 'Henan Topeco Clean Import & Export Co., Ltd.'
    '28' = 'Anqing Jiaxin Medical Technology Co., Ltd.'
    '29' = 'Linghai Zhanwang Biotechnology Co., Ltd.'
    '30' = 'Shanghai Chuangsi Youxiang Trading Co., Ltd.'
    '31' = 'Changshu Power Clean Co., Ltd.'
    '32' = 'Henan Zenpe Technology Co., Ltd.'
    '33' = 'Zhejiang Furuisen Spunlaced Non-Wovens Co., Ltd.'
    '34' = 'Fuyang Yangyang Health Technology Co., Ltd.'
    '35' = 'Dongguan Ruisheng Textile Co., Ltd.'
    '36' = 'Zhejiang Qimei Commodity Co., Ltd.'
    '37' = 'Yiwu Chen Ju Electronic Commerce Co., Ltd.'
    '38' = 'Colorful medium inc.'
    '39' = 'King Tech Co., Ltd.'
    '40' = 'Zhejiang Huashun Technology Co., Ltd.'
    '41' = 'Hangzhou Shengbo Cleaning Product Co., Ltd.'
    '42' = 'Ningbo Yachen Import & Export Co., Ltd.'
    '43' = 'Shanghai Kingmax Commodity Co., Ltd.'
    '44' = 'Xiamen Mk Health Care Product Co., Ltd.'
    '45' = 'Changshu He Gui Textiles Co., Ltd.'
    '46' = 'Pujiang Hailan Garments Co., Ltd.'
    '47' = 'Dongguan Winall Paper Co., Ltd.'
    '48' = 'Changshu Beisimei Knitting Textile Co., Ltd.'
    '49' = 'Shanghai Mtg Cleaning Material Co., Ltd.'
}

foreach ($row in $companies.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $companies[$row]
}
